# fix subset method definition, remove unneeded plotSpec() override
#
# Adds a new "components" worksheet (status-tracking table for the
# `components`-class methods, mirroring the existing fGroups/mslists/
# formulas/compounds sheets) after the "compounds" sheet, makes it the
# active/selected sheet, and resets the previously-active "compounds"
# sheet's selection back to its header row.

$wb = $excel.ActiveWorkbook

# --- add the new sheet after the last (compounds) sheet ------------------
$compounds = $wb.Worksheets.Item("compounds")
$newSheet = $wb.Worksheets.Add($null, $compounds)
$newSheet.Name = "components"

# --- header row ------------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "as-is"
$newSheet.Cells.Item(1, 3).Value = "almost as-is"
$newSheet.Cells.Item(1, 4).Value = "implement"
$newSheet.Cells.Item(1, 5).Value = "not supported"
$newSheet.Cells.Item(1, 6).Value = "ionize"
$newSheet.Cells.Item(1, 7).Value = "done"

# --- status rows -------------------------------------------------------
# row, A, B, C, D, E, G, H
$rows = @(
    @(2,  "$",              "X", $null, $null, $null, "X", $null),
    @(3,  "[",               $null, "X", $null, $null, "X", $null),
    @(4,  "[[",              "X", $null, $null, $null, "X", $null),
    @(5,  "as.data.table",   "X", $null, $null, $null, "X", $null),
    @(6,  "componentInfo",   "X", $null, $null, $null, "X", $null),
    @(7,  "componentTable",  "X", $null, $null, $null, "X", $null),
    @(8,  "consensus",       $null, $null, $null, "X", $null, $null),
    @(9,  "filter",          $null, "X", $null, $null, "X", $null),
    @(10, "findFGroup",      "X", $null, $null, $null, "X", $null),
    @(11, "groupNames",      "X", $null, $null, $null, "X", $null),
    @(12, "initialize",      $null, "X", $null, $null, "X", $null),
    @(13, "length",          "X", $null, $null, $null, "X", $null),
    @(14, "names",           "X", $null, $null, $null, "X", $null),
    @(15, "plotEIC",         "X", $null, "X", $null, "X", "Seems enough, assuming we're not planning to merge components"),
    @(16, "plotEICHash",     "X", $null, $null, $null, "X", $null),
    @(17, "plotSpec",        "X", $null, "X", $null, "X", "Seems enough, assuming we're not planning to merge components"),
    @(18, "plotSpecHash",    "X", $null, $null, $null, "X", $null),
    @(19, "show",            $null, "X", $null, $null, "X", $null)
)

foreach ($row in $rows) {
    $r = $row[0]
    $newSheet.Cells.Item($r, 1).Value = $row[1]
    if ($row[2]) { $newSheet.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3]) { $newSheet.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4]) { $newSheet.Cells.Item($r, 4).Value = $row[4] }
    if ($row[5]) { $newSheet.Cells.Item($r, 5).Value = $row[5] }
    if ($row[6]) { $newSheet.Cells.Item($r, 7).Value = $row[6] }
    if ($row[7]) { $newSheet.Cells.Item($r, 8).Value = $row[7] }
}

$newSheet.Columns.Item(1).ColumnWidth = 16.140625

# --- selection / activation ----------------------------------------------
$newSheet.Activate()
$newSheet.Range("H17").Select()

$compounds.Activate()
$compounds.Range("B1:G1").Select()

$newSheet.Activate()
